$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 44684
$ws.Range("J108").Value = 44684
$ws.Range("L108").Value = 44684
$ws.Range("N108").Value = -52364
$ws.Range("H117").Value = 48692
$ws.Range("J117").Value = 48692
$ws.Range("L117").Value = 48692
$ws.Range("N117").Value = -57870
$ws.Range("H120").Value = 48311.6
$ws.Range("J120").Value = 48311.6
$ws.Range("L120").Value = 48311.6
$ws.Range("N120").Value = -57987.6
$ws.Range("H128").Value = 41575
$ws.Range("J128").Value = 41575
$ws.Range("L128").Value = 41575
$ws.Range("N128").Value = -51535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22980.953
$ws.Range("I32").Value = 23485.848
$ws.Range("J32").Value = 16333.167
$ws.Range("K32").Value = 23485.848
$ws.Range("L32").Value = 16333.167
$ws.Range("M32").Value = -23198.848
$ws.Range("N32").Value = -16907.167
$ws.Range("H107").Value = 42492
$ws.Range("J107").Value = 42492
$ws.Range("L107").Value = 42492
$ws.Range("N107").Value = -50172
$ws.Range("H109").Value = 41877
$ws.Range("J109").Value = 41877
$ws.Range("L109").Value = 41877
$ws.Range("N109").Value = -44651
$ws.Range("H118").Value = 44946
$ws.Range("J118").Value = 44946
$ws.Range("L118").Value = 44946
$ws.Range("N118").Value = -48260
$ws.Range("H125").Value = 49992
$ws.Range("J125").Value = 49992
$ws.Range("L125").Value = 49992
$ws.Range("N125").Value = -59832
$ws.Range("H128").Value = 50421
$ws.Range("J128").Value = 50421
$ws.Range("L128").Value = 50421
$ws.Range("N128").Value = -60381
$ws.Range("H130").Value = 42432
$ws.Range("J130").Value = 42432
$ws.Range("L130").Value = 42432
$ws.Range("N130").Value = -52472
$ws.Range("H131").Value = 50661
$ws.Range("J131").Value = 50661
$ws.Range("L131").Value = 50661
$ws.Range("N131").Value = -60741

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2181.6365
$ws.Range("I105").Value = 1942.7778
$ws.Range("J105").Value = 2347
$ws.Range("K105").Value = 1942.7778
$ws.Range("L105").Value = 2347
$ws.Range("M105").Value = -195.7778000000001
$ws.Range("N105").Value = -5841
$ws.Range("H107").Value = 2227.84
$ws.Range("I107").Value = 2033
$ws.Range("J107").Value = 3250.75
$ws.Range("K107").Value = 2033
$ws.Range("L107").Value = 3250.75
$ws.Range("M107").Value = -113
$ws.Range("N107").Value = -7090.75
$ws.Range("H119").Value = 40507.332
$ws.Range("J119").Value = 40507.332
$ws.Range("L119").Value = 40507.332
$ws.Range("N119").Value = -50183.332
$ws.Range("H120").Value = 44727.668
$ws.Range("J120").Value = 44727.668
$ws.Range("L120").Value = 44727.668
$ws.Range("N120").Value = -54403.668
$ws.Range("H125").Value = 48992
$ws.Range("J125").Value = 48992
$ws.Range("L125").Value = 48992
$ws.Range("N125").Value = -58832
$ws.Range("H126").Value = 49884
$ws.Range("J126").Value = 49884
$ws.Range("L126").Value = 49884
$ws.Range("N126").Value = -59764

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49379.75
$ws.Range("J20").Value = 49379.75
$ws.Range("L20").Value = 49379.75
$ws.Range("N20").Value = -49851.75
$ws.Range("H22").Value = 346.375
$ws.Range("I22").Value = 428.66666
$ws.Range("J22").Value = 99.5
$ws.Range("K22").Value = 428.66666
$ws.Range("L22").Value = 99.5
$ws.Range("M22").Value = -78.66665999999998
$ws.Range("N22").Value = -799.5
$ws.Range("H30").Value = 49379.75
$ws.Range("J30").Value = 49379.75
$ws.Range("L30").Value = 49379.75
$ws.Range("N30").Value = -49561.75
$ws.Range("H116").Value = 46110.5
$ws.Range("J116").Value = 46110.5
$ws.Range("L116").Value = 46110.5
$ws.Range("N116").Value = -55288.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H128").Value = 49379.75
$ws.Range("J128").Value = 49379.75
$ws.Range("L128").Value = 49379.75
$ws.Range("N128").Value = -59339.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 924.03
$ws.Range("J131").Value = 929.9382000000001
$ws.Range("L131").Value = 2789.8146
$ws.Range("N131").Value = -12869.8146
$ws.Range("H132").Value = 1910
$ws.Range("I132").Value = 1400
$ws.Range("J132").Value = 1966.6666
$ws.Range("K132").Value = 12600
$ws.Range("L132").Value = 17699.9994
$ws.Range("M132").Value = -10070
$ws.Range("N132").Value = -22759.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 33800
$ws.Range("J110").Value = 33800
$ws.Range("L110").Value = 33800
$ws.Range("N110").Value = -41980
$ws.Range("H130").Value = 50661.332
$ws.Range("J130").Value = 50661.332
$ws.Range("L130").Value = 50661.332
$ws.Range("N130").Value = -60701.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2168.8096
$ws.Range("I7").Value = 1874.2858
$ws.Range("J7").Value = 2757.8572
$ws.Range("K7").Value = 1874.2858
$ws.Range("L7").Value = 2757.8572
$ws.Range("M7").Value = -1762.2858
$ws.Range("N7").Value = -2981.8572
$ws.Range("H111").Value = 43626.332
$ws.Range("J111").Value = 43626.332
$ws.Range("L111").Value = 43626.332
$ws.Range("N111").Value = -51806.332
$ws.Range("H121").Value = 44137.5
$ws.Range("J121").Value = 44137.5
$ws.Range("L121").Value = 44137.5
$ws.Range("N121").Value = -47631.5
$ws.Range("H126").Value = 2168.8096
$ws.Range("I126").Value = 1874.2858
$ws.Range("J126").Value = 2757.8572
$ws.Range("K126").Value = 5622.857400000001
$ws.Range("L126").Value = 8273.571599999999
$ws.Range("M126").Value = -3152.857400000001
$ws.Range("N126").Value = -13213.5716
$ws.Range("H130").Value = 37521.555
$ws.Range("J130").Value = 37521.555
$ws.Range("L130").Value = 37521.555
$ws.Range("N130").Value = -47561.555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 39221.25
$ws.Range("J119").Value = 39221.25
$ws.Range("L119").Value = 39221.25
$ws.Range("N119").Value = -48897.25
$ws.Range("H120").Value = 43966.4
$ws.Range("J120").Value = 43966.4
$ws.Range("L120").Value = 43966.4
$ws.Range("N120").Value = -53642.4
